$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K11:K16").Value = "Passed"

$ws.Range("I2").FormulaArray = "=SUMPRODUCT((K10:K16=""Passed"")+(K10:K16=""Failed""))"
$ws.Range("I4").Formula = "=COUNTIF(K10:K16, ""Passed"")"
